$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sequential Tournament IDs (column K) for rows 2..392, starting at 100.
# Each distinct original Tournament ID is replaced by the next sequential
# integer (100, 101, 102, ...) assigned in order of first appearance.
$newIds = @(
    100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 100, 110, 100, 111, 112, 113, 114, 115, 116, 113, 117, 118, 119, 120, 100,
    121, 122, 100, 123, 124, 113, 125, 126, 113, 127, 128, 129, 100, 124, 130, 100, 131, 124, 113, 132, 133, 134, 113, 135, 136,
    137, 100, 138, 124, 139, 100, 140, 124, 113, 141, 142, 143, 113, 144, 145, 146, 100, 124, 147, 124, 148, 113, 149, 150, 100,
    151, 113, 152, 153, 154, 124, 155, 156, 157, 124, 158, 113, 159, 160, 156, 161, 162, 156, 163, 164, 165, 166, 167, 156, 168,
    169, 170, 171, 156, 172, 173, 174, 175, 176, 156, 177, 178, 179, 180, 181, 182, 156, 183, 184, 185, 186, 187, 188, 189, 190,
    191, 192, 193, 194, 195, 196, 197, 198, 199, 200, 201, 202, 203, 204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215,
    205, 216, 217, 218, 219, 220, 221, 222, 223, 224, 205, 225, 226, 227, 228, 229, 230, 231, 232, 205, 233, 234, 235, 236, 237,
    238, 239, 240, 241, 242, 243, 244, 245, 246, 247, 248, 249, 250, 251, 252, 253, 251, 252, 254, 255, 251, 252, 256, 257, 258,
    251, 252, 259, 260, 261, 251, 252, 262, 263, 251, 252, 264, 265, 251, 252, 266, 251, 267, 268, 269, 270, 271, 272, 252, 273,
    274, 275, 276, 277, 278, 279, 280, 281, 282, 283, 284, 285, 286, 287, 288, 289, 290, 291, 292, 293, 294, 295, 296, 297, 298,
    299, 300, 301, 299, 302, 303, 299, 304, 305, 299, 306, 307, 299, 308, 309, 299, 310, 311, 312, 299, 313, 314, 315, 316, 317,
    318, 313, 317, 319, 320, 317, 321, 322, 317, 323, 324, 317, 325, 317, 326, 327, 328, 317, 329, 330, 331, 332, 333, 334, 335,
    336, 334, 337, 338, 334, 339, 340, 334, 341, 342, 334, 343, 344, 334, 345, 346, 347, 334, 348, 349, 350, 351, 352, 353, 348,
    354, 355, 356, 357, 358, 359, 360, 361, 362, 363, 364, 365, 366, 362, 367, 368, 369, 370, 371, 372, 373, 374, 375, 376, 377,
    378, 379, 380, 381, 382, 383, 378, 384, 385, 386, 387, 388, 389, 378, 390, 391, 392, 393, 394, 395, 396, 397, 398, 378, 399,
    391, 400, 401, 393, 392, 402, 403, 404, 405, 378, 406, 407, 400, 408, 392, 409
)

$startRow = 2
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 11).Value = $newIds[$i]
}
